{"js": "const replacements = [\n  [\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"],\n  [\"30\u00d758=\", \"33\u00d771=\"],\n  [\"17\u00d727=\", \"57\u00d744=\"],\n  [\"11\u00d747=\", \"27\u00d739=\"],\n  [\"26\u00d798=\", \"37\u00d767=\"],\n  [\"57\u00d781=\", \"50\u00d717=\"],\n  [\"54\u00d775=\", \"18\u00d742=\"],\n  [\"74\u00d792=\", \"22\u00d731=\"],\n  [\"66\u00d798=\", \"64\u00d712=\"],\n  [\"98\u00d729=\", \"14\u00d797=\"],\n  [\"25\u00d781=\", \"75\u00d730=\"],\n  [\"34\u00d775=\", \"75\u00d760=\"],\n  [\"35\u00d786=\", \"70\u00d715=\"],\n  [\"95\u00d776=\", \"29\u00d756=\"],\n  [\"45\u00d797=\", \"60\u00d746=\"],\n  [\"98\u00d798=\", \"62\u00d744=\"],\n  [\"40\u00d769=\", \"22\u00d745=\"],\n  [\"99\u00d721=\", \"88\u00d738=\"],\n  [\"16\u00d739=\", \"36\u00d798=\"],\n  [\"58\u00d725=\", \"62\u00d732=\"],\n  [\"54\u00d747=\", \"41\u00d784=\"],\n  [\"29\u00d755=\", \"72\u00d739=\"],\n  [\"54\u00d758=\", \"28\u00d740=\"],\n  [\"40\u00d736=\", \"45\u00d759=\"],\n  [\"93\u00d714=\", \"67\u00d765=\"],\n  [\"53\u00d754=\", \"87\u00d722=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"),\n    @(\"30\u00d758=\", \"33\u00d771=\"),\n    @(\"17\u00d727=\", \"57\u00d744=\"),\n    @(\"11\u00d747=\", \"27\u00d739=\"),\n    @(\"26\u00d798=\", \"37\u00d767=\"),\n    @(\"57\u00d781=\", \"50\u00d717=\"),\n    @(\"54\u00d775=\", \"18\u00d742=\"),\n    @(\"74\u00d792=\", \"22\u00d731=\"),\n    @(\"66\u00d798=\", \"64\u00d712=\"),\n    @(\"98\u00d729=\", \"14\u00d797=\"),\n    @(\"25\u00d781=\", \"75\u00d730=\"),\n    @(\"34\u00d775=\", \"75\u00d760=\"),\n    @(\"35\u00d786=\", \"70\u00d715=\"),\n    @(\"95\u00d776=\", \"29\u00d756=\"),\n    @(\"45\u00d797=\", \"60\u00d746=\"),\n    @(\"98\u00d798=\", \"62\u00d744=\"),\n    @(\"40\u00d769=\", \"22\u00d745=\"),\n    @(\"99\u00d721=\", \"88\u00d738=\"),\n    @(\"16\u00d739=\", \"36\u00d798=\"),\n    @(\"58\u00d725=\", \"62\u00d732=\"),\n    @(\"54\u00d747=\", \"41\u00d784=\"),\n    @(\"29\u00d755=\", \"72\u00d739=\"),\n    @(\"54\u00d758=\", \"28\u00d740=\"),\n    @(\"40\u00d736=\", \"45\u00d759=\"),\n    @(\"93\u00d714=\", \"67\u00d765=\"),\n    @(\"53\u00d754=\", \"87\u00d722=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    [void]$find.Execute(\n        $old,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $new,\n        2\n    )\n}\n"}
